# Updated symbol list (coinranking scrape refresh) - updates Price (col D)
# and Volume(1h) (col E) text values for the rows whose figures moved.
#
# The source cells are plain text (not numbers) even though most of the
# strings look numeric/percent-like. A bare `.Value = "309.13"` would get
# auto-coerced by Excel into a real number (and "0.87%" into a fraction),
# which would change the cell's type/formatting and not match the sheet's
# original text-cell layout. Prefixing with a leading apostrophe forces
# Excel to store the literal text, and the follow-up ClearFormats() strips
# the "quote prefix" cell style Excel applies for that trick, so the cell
# ends up as plain text with no style change versus the original - just
# like the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    $rng.Value = "'" + $Text
    $rng.ClearFormats()
}

# Row 2 - BNB
Set-TextValue "D2" "309.13"
Set-TextValue "E2" "0.87%"

# Row 3 - OKB
Set-TextValue "D3" "40.97"
Set-TextValue "E3" "1.48%"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.125"

# Row 5 - Cronos
Set-TextValue "D5" "0.07622"
Set-TextValue "E5" "0.43%"

# Row 6 - GateToken
Set-TextValue "D6" "4.282"
Set-TextValue "E6" "0.74%"

# Row 7 - FTXToken
Set-TextValue "D7" "1.618"
Set-TextValue "E7" "1.97%"

# Row 8 - BTSEToken
Set-TextValue "D8" "2.488"
Set-TextValue "E8" "1.79%"

# Row 9 - MXToken
Set-TextValue "D9" "0.9102"
Set-TextValue "E9" "0.67%"

# Row 10 - LiechtensteinCryptoassetsExchange
Set-TextValue "D10" "0.1185"
Set-TextValue "E10" "18.92%"

# Row 11 - WazirX
Set-TextValue "D11" "0.1821"
Set-TextValue "E11" "3.30%"

# Row 12 - MandalaExchangeToken
Set-TextValue "D12" "0.09112"
Set-TextValue "E12" "1.19%"

# Row 13 - BitrueCoin
Set-TextValue "D13" "0.04275"
Set-TextValue "E13" "-2.69%"

# Row 14 - BitMartToken
Set-TextValue "E14" "-0.77%"

# Row 15 - BitForexToken
Set-TextValue "D15" "0.001261"
Set-TextValue "E15" "0.37%"

# Row 16 - TigerCash
Set-TextValue "D16" "0.005781"
Set-TextValue "E16" "-0.77%"

# Row 17 - LEO
Set-TextValue "D17" "3.351"
Set-TextValue "E17" "-0.54%"

# Row 18 - BitpandaEcosystemToken
Set-TextValue "E18" "-0.64%"

# Row 19 - MCDex
Set-TextValue "D19" "6.912"
Set-TextValue "E19" "1.03%"

# Row 20 - ProBitToken
Set-TextValue "E20" "4.40%"

# Row 22 - CoinExToken
Set-TextValue "D22" "0.04037"
Set-TextValue "E22" "-3.01%"

# Row 23 - BitKan
Set-TextValue "E23" "4.36%"

# Row 24 - HotbitToken
Set-TextValue "D24" "0.004120"
Set-TextValue "E24" "1.31%"

# Row 25 - NitroEx
Set-TextValue "E25" "-2.47%"

# Row 26 - UpBots
Set-TextValue "E26" "24.33%"

# Row 38 - One
Set-TextValue "D38" "0.02420"
Set-TextValue "E38" "1.20%"

# Row 39 - IDEX
Set-TextValue "D39" "0.05238"
Set-TextValue "E39" "2.28%"

# Row 40 - KickToken
Set-TextValue "D40" "0.007787"
Set-TextValue "E40" "-0.85%"

# Row 41 - BKEXToken
Set-TextValue "E41" "-0.22%"

# Row 42 - Dexo
Set-TextValue "D42" "0.006806"
Set-TextValue "E42" "-4.13%"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.008069"
Set-TextValue "E44" "-3.75%"

# Row 45 - PooCoin
Set-TextValue "D45" "0.3073"
Set-TextValue "E45" "-7.40%"

# Row 46 - CoinLion
Set-TextValue "D46" "0.00006900"
Set-TextValue "E46" "6.80%"

# Row 47 - Kangarootoken
Set-TextValue "E47" "-0.16%"

# Row 48 - BOLO
Set-TextValue "D48" "0.09853"
Set-TextValue "E48" "1,624.67%"

# Row 50 - CryptobidCoin
Set-TextValue "D50" "0.00002103"
Set-TextValue "E50" "-0.16%"

# Row 51 - SpecialPowerGold
Set-TextValue "E51" "-0.16%"
